# "first review of the entire workflow"
# Updated the lower/upper RLIe confidence-interval bounds (columns F and G)
# for the bootstrap re-run across GET_L2 / GET_L3 / Aggregate rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

# GET_L2 rows
$ws.Range("G5").Value  = 0.78
$ws.Range("F6").Value  = 0.8363636363636364
$ws.Range("F14").Value = 0.5834586466165413
$ws.Range("G14").Value = 0.7097744360902256
$ws.Range("F16").Value = 0.5714285714285714
$ws.Range("F17").Value = 0.5714285714285714
$ws.Range("G17").Value = 0.6947368421052631
$ws.Range("F18").Value = 0.8333333333333334
$ws.Range("G18").Value = 0.9055555555555556
$ws.Range("F19").Value = 0.8044444444444444
$ws.Range("G19").Value = 0.8822222222222222
$ws.Range("F20").Value = 0.7977777777777778
$ws.Range("F21").Value = 0.7944444444444445
$ws.Range("G21").Value = 0.8733333333333333
$ws.Range("F22").Value = 0.8934782608695653
$ws.Range("G23").Value = 0.9652173913043478
$ws.Range("F24").Value = 0.8695652173913043
$ws.Range("G24").Value = 0.9673913043478261
$ws.Range("F25").Value = 0.8673913043478261

# GET_L3 rows
$ws.Range("G27").Value = 0.8
$ws.Range("F30").Value = 0.8272727272727273
$ws.Range("G42").Value = 0.7075757575757575
$ws.Range("G43").Value = 0.6909469696969691
$ws.Range("G44").Value = 0.6924242424242424
$ws.Range("G45").Value = 0.6924242424242424
$ws.Range("F46").Value = 0.8842105263157894
$ws.Range("F47").Value = 0.8842105263157894
$ws.Range("F48").Value = 0.8578947368421053
$ws.Range("F49").Value = 0.8630263157894738
$ws.Range("G50").Value = 0.9105223880597004
$ws.Range("F51").Value = 0.7432835820895523
$ws.Range("G52").Value = 0.8746268656716418
$ws.Range("G53").Value = 0.8746268656716418
$ws.Range("G55").Value = 0.8826666666666667
$ws.Range("F57").Value = 0.7466666666666666
$ws.Range("G62").Value = 0.9746031746031746
$ws.Range("F64").Value = 0.8444444444444444

# Aggregate rows
$ws.Range("F70").Value = 0.7935205183585313
$ws.Range("G70").Value = 0.8488120950323974
$ws.Range("F71").Value = 0.7688876889848812
$ws.Range("G71").Value = 0.8267926565874728
$ws.Range("G72").Value = 0.8241900647948164
$ws.Range("F73").Value = 0.7645788336933046
